$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from A464 (date column) to new date rows A465:A491
$ws.Cells.Item(464, 1).Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(465, 44539, 4, 21, 164.8998822143699),
    @(466, 44540, 8, 29, 227.7188849627012),
    @(467, 44541, 0, 28, 219.8665096191598),
    @(468, 44542, 4, 26, 204.161758932077),
    @(469, 44543, 7, 30, 235.5712603062426),
    @(470, 44544, 3, 27, 212.0141342756183),
    @(471, 44545, 0, 26, 204.161758932077),
    @(472, 44546, 5, 27, 212.0141342756183),
    @(473, 44547, 13, 32, 251.2760109933255),
    @(474, 44548, 5, 37, 290.5378877110326),
    @(475, 44550, 6, 39, 306.2426383981154),
    @(476, 44551, 6, 38, 298.390263054574),
    @(477, 44552, 1, 36, 282.6855123674911),
    @(478, 44553, 9, 45, 353.356890459364),
    @(479, 44554, 10, 50, 392.618767177071),
    @(480, 44555, 5, 42, 329.7997644287397),
    @(481, 44556, 3, 40, 314.0950137416569),
    @(482, 44557, 8, 42, 329.7997644287397),
    @(483, 44558, 1, 37, 290.5378877110326),
    @(484, 44559, 5, 41, 321.9473890851983),
    @(485, 44560, 23, 55, 431.8806438947781),
    @(486, 44561, 17, 62, 486.8472712995681),
    @(487, 44562, 19, 76, 596.780526109148),
    @(488, 44563, 18, 91, 714.5661562622694),
    @(489, 44564, 11, 94, 738.1232822928936),
    @(490, 44565, 2, 95, 745.975657636435),
    @(491, 44566, 17, 107, 840.2041617589321)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Host "Update complete"
